# Trade #13 closed at 2026-02-16 21:54:52 - leadlag DOWN +0.000%
#
# Appends the new trade row (Trade #13) to both the "All Trades" sheet
# (next empty row 14) and the strategy-specific "leadlag" sheet (next
# empty row 13), mirroring the existing rows already present there.

$wb = $excel.ActiveWorkbook

function Add-TradeRow {
    param(
        [object]$ws,
        [int]$row
    )

    # Column A: Trade # (number)
    $ws.Cells.Item($row, 1).Value = 13

    # Column B: Date - force Text format so Excel does not reinterpret
    # the "2026-02-16" string as a date value (matches the existing
    # rows, which are stored as plain text).
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"

    # Column C: Time (text)
    $ws.Cells.Item($row, 3).Value = "21:54:52"

    # Column D: Strategy
    $ws.Cells.Item($row, 4).Value = "leadlag"

    # Column E: Side
    $ws.Cells.Item($row, 5).Value = "DOWN"

    # Column F: Entry Price
    $ws.Cells.Item($row, 6).Value = 68243.81

    # Column G: Exit Price (blank - trade still OPEN)

    # Column H: Status
    $ws.Cells.Item($row, 8).Value = "OPEN"

    # Column I: P&L %
    $ws.Cells.Item($row, 9).Value = 0

    # Column J: P&L $
    $ws.Cells.Item($row, 10).Value = 0

    # Column K: Capital After
    $ws.Cells.Item($row, 11).Value = 100

    # Column L: Confidence
    $ws.Cells.Item($row, 12).Value = 0.75

    # Column M: Entry Reason
    $ws.Cells.Item($row, 13).Value = "Coinbase leading with -0.100% move"

    # Column N: Exit Reason (blank - trade still OPEN)

    # Column O: Duration (min)
    $ws.Cells.Item($row, 15).Value = 0
}

# "All Trades" sheet gains the new trade as row 14
$wsAll = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAll 14

# "leadlag" strategy sheet gains the same new trade as row 13
$wsLeadLag = $wb.Worksheets.Item("leadlag")
Add-TradeRow $wsLeadLag 13
